# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions" — updates Price (D) /
# Volume(1h) (E) values for every coin row, and fixes four coin rows whose
# Coin/Link (B/C) pairs were swapped (rows 14/15, 19/20, 28/29, 42/43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.026.74"
$ws.Range("E2").Value = "  +2.98%  "

$ws.Range("D3").Value = "1.698.56"
$ws.Range("E3").Value = "  +2.46%  "

$ws.Range("D4").Value = "'0.9948"
$ws.Range("E4").Value = "  -0.89%  "

$ws.Range("D5").Value = "'229.54"
$ws.Range("E5").Value = "  +4.61%  "

$ws.Range("D6").Value = "'0.5407"
$ws.Range("E6").Value = "  +3.11%  "

$ws.Range("D7").Value = "'0.9952"
$ws.Range("E7").Value = "  -0.92%  "

$ws.Range("D8").Value = "'0.2712"
$ws.Range("E8").Value = "  +2.16%  "

$ws.Range("D9").Value = "'0.06486"
$ws.Range("E9").Value = "  +2.48%  "

$ws.Range("D10").Value = "'21.44"
$ws.Range("E10").Value = "  +3.49%  "

$ws.Range("D11").Value = "'0.07734"
$ws.Range("E11").Value = "  -0.57%  "

$ws.Range("D12").Value = "'4.722"
$ws.Range("E12").Value = "  +3.50%  "

$ws.Range("D13").Value = "1.685.48"
$ws.Range("E13").Value = "  +2.54%  "

$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.5977"
$ws.Range("E14").Value = "  +5.58%  "

$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "1.906.24"
$ws.Range("E15").Value = "  +1.13%  "

$ws.Range("D16").Value = "0.0₅8288"
$ws.Range("E16").Value = "  +2.29%  "

$ws.Range("D17").Value = "'67.37"
$ws.Range("E17").Value = "  +2.94%  "

$ws.Range("D18").Value = "26.850.66"
$ws.Range("E18").Value = "  +2.35%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'4.788"
$ws.Range("E19").Value = "  +1.12%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'0.9984"
$ws.Range("E20").Value = "  -0.61%  "

$ws.Range("D21").Value = "'209.14"
$ws.Range("E21").Value = "  +8.03%  "

$ws.Range("D22").Value = "'10.94"
$ws.Range("E22").Value = "  +6.17%  "

$ws.Range("D23").Value = "'6.150"
$ws.Range("E23").Value = "  +1.88%  "

$ws.Range("D24").Value = "'0.9958"
$ws.Range("E24").Value = "  -0.89%  "

$ws.Range("D25").Value = "'146.65"
$ws.Range("E25").Value = "  +1.81%  "

$ws.Range("D26").Value = "'0.1238"
$ws.Range("E26").Value = "  +2.78%  "

$ws.Range("D27").Value = "'7.395"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'16.55"
$ws.Range("E28").Value = "  +3.28%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'1.571"
$ws.Range("E29").Value = "  +4.28%  "

$ws.Range("D30").Value = "'0.05696"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("D31").Value = "'1.319"
$ws.Range("E31").Value = "  +3.36%  "

$ws.Range("D32").Value = "'3.616"
$ws.Range("E32").Value = "  +3.35%  "

$ws.Range("D33").Value = "'3.508"
$ws.Range("E33").Value = "  +3.52%  "

$ws.Range("D34").Value = "'1.630"
$ws.Range("E34").Value = "  +2.18%  "

$ws.Range("D35").Value = "'0.9820"
$ws.Range("E35").Value = "  +3.72%  "

$ws.Range("D36").Value = "'2.824"
$ws.Range("E36").Value = "  +0.76%  "

$ws.Range("D37").Value = "'2.408"
$ws.Range("E37").Value = "  +0.30%  "

$ws.Range("D38").Value = "'0.5832"
$ws.Range("E38").Value = "  +1.26%  "

$ws.Range("D39").Value = "'0.01629"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").Value = "'5.975"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("D41").Value = "1.058.63"
$ws.Range("E41").Value = "  +1.82%  "

$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.8385"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9956"
$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").Value = "'102.69"
$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").Value = "1.822.74"
$ws.Range("E45").Value = "  +1.53%  "

$ws.Range("D46").Value = "'59.62"
$ws.Range("E46").Value = "  +2.10%  "

$ws.Range("D48").Value = "'0.9956"
$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("D49").Value = "'8.035"
$ws.Range("E49").Value = "  +0.44%  "

$ws.Range("D50").Value = "'0.4324"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("D51").Value = "'0.05216"
$ws.Range("E51").Value = "  -1.96%  "
